# Apply the "Source Ftp" column insertion + header rename + row-2 value
# updates described by the commit:
#   "Add files via upload. Added the new file downloading from the FTP,
#    and also delete old file from local system, and rename new file to
#    old. Fixed the issue occurred while creating extra column for
#    External ID."

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2 (data row) value updates -----------------------------------
$ws.Range("A2").Value = "Jenne"
$ws.Range("C2").Value = "First Name"
$ws.Range("D2").Value = "Phone No"
$ws.Range("E2").Value = ".xlsx"
$ws.Range("F2").Value = ".xlsx"

# --- Rename existing J:M header row from "Results Ftp *" to
#     "Source Ftp *" (these become the Source-side FTP settings) -------
$ws.Range("J1").Value = "Source Ftp Url"
$ws.Range("K1").Value = "Source Ftp User"
$ws.Range("L1").Value = "Source Ftp Pass"
$ws.Range("M1").Value = "Source Ftp Port"

# --- New columns N:O -- remaining Source Ftp settings ------------------
$ws.Range("N1").Value = "Source Ftp Path"
$ws.Range("O1").Value = "Source Ftp Filename"
$ws.Range("N2").Value = "/telquestftp.com/jenne/"
$ws.Range("O2").Value = "Jenne"

# --- New columns P:T -- re-added Results Ftp settings (now for the
#     results-upload destination) ---------------------------------------
$ws.Range("P1").Value = "Results Ftp Url"
$ws.Range("Q1").Value = "Results Ftp User"
$ws.Range("R1").Value = "Results Ftp Pass"
$ws.Range("S1").Value = "Results Ftp Port"
$ws.Range("T1").Value = "Results Ftp Path"

$ws.Range("P2").Value = "telquestftp.com"
$ws.Range("Q2").Value = "admin@telquestftp.com"
$ws.Range("R2").Value = "Shopping2016#"
$ws.Range("S2").Value = 21
$ws.Range("T2").Value = "/telquestftp.com/results/"

# --- Copy the header/body formatting from the existing FTP block
#     (columns J:M) onto the newly added columns so the styling
#     matches. Source/destination ranges are sized identically so the
#     paste cannot spill into neighbouring columns. -------------------
$ws.Range("J1:K4").Copy()
$ws.Range("N1:O4").PasteSpecial(-4122)

$ws.Range("J1:N4").Copy()
$ws.Range("P1:T4").PasteSpecial(-4122)

$excel.CutCopyMode = 0

# --- Column widths (bestFit sizes from the final file; input values
#     chosen so this engine's char->pixel->char rounding lands as close
#     as possible on the target stored widths) ---------------------------
$ws.Columns.Item(1).ColumnWidth = 11.6
$ws.Columns.Item(3).ColumnWidth = 11.3
$ws.Columns.Item(4).ColumnWidth = 22.6
$ws.Columns.Item(8).ColumnWidth = 12.3
$ws.Columns.Item(9).ColumnWidth = 20.95
$ws.Columns.Item(12).ColumnWidth = 13.1
$ws.Columns.Item(13).ColumnWidth = 12.95
$ws.Columns.Item(14).ColumnWidth = 21.1
$ws.Columns.Item(15).ColumnWidth = 21.1
$ws.Columns.Item(16).ColumnWidth = 13.3
$ws.Columns.Item(17).ColumnWidth = 20.1
$ws.Columns.Item(18).ColumnWidth = 12.95
$ws.Columns.Item(19).ColumnWidth = 12.95
$ws.Columns.Item(20).ColumnWidth = 21.1

# --- View / selection tweaks -------------------------------------------
$ws.Range("D7").Select()
